$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.573.08"
$ws.Range("E2").Value = "  -1.30%  "

$ws.Range("D3").Value = "3.012.50"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.02"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.84%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  -2.53%  "

$ws.Range("D9").Value = "3.011.90"
$ws.Range("E9").Value = "  -1.68%  "

$ws.Range("E10").Value = "  -4.05%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.81"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.31%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.459"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.71%  "

$ws.Range("E13").Value = "  -2.96%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.68%  "

$ws.Range("E15").Value = "  +2.06%  "

$ws.Range("D16").Value = "3.507.16"
$ws.Range("E16").Value = "  -1.71%  "

$ws.Range("E17").Value = "  -1.16%  "

$ws.Range("D18").Value = "62.516.17"
$ws.Range("E18").Value = "  -1.29%  "

$ws.Range("D19").Value = "3.010.26"
$ws.Range("E19").Value = "  -1.76%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "459.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.56%  "

$ws.Range("E21").Value = "  -2.74%  "

$ws.Range("E22").Value = "  -2.63%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.43"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.61%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.53%  "

$ws.Range("B25").Value = "Fetch.AI"
$ws.Range("C25").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.62%  "

$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.93%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.07%  "

$ws.Range("E28").Value = "  +0.10%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.06%  "

$ws.Range("E30").Value = "  -2.57%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.49%  "

$ws.Range("E32").Value = "  -5.86%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "28.13"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.58%  "

$ws.Range("E34").Value = "  -2.31%  "

$ws.Range("D35").Value = "0.0₃0813"
$ws.Range("E35").Value = "  -1.46%  "

$ws.Range("E36").Value = "  -3.45%  "

$ws.Range("E37").Value = "  -3.56%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.13"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.25%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.41"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.39%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.14"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.50%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -13.09%  "

$ws.Range("E42").Value = "  +4.28%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "390.43"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -11.15%  "

$ws.Range("E44").Value = "  -1.68%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.270"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.58%  "

$ws.Range("D46").Value = "2.737.08"
$ws.Range("E46").Value = "  -4.01%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "37.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.79%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.52"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.66%  "

$ws.Range("E50").Value = "  -1.53%  "

$ws.Range("E51").Value = "  -0.79%  "
